$d = $word.ActiveDocument

# Locate the paragraph that contains the sentence about the access token,
# searching by a stable substring so we don't depend on a hard-coded index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*will receive a token that will grant you*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the target paragraph"
}

$pStart = $target.Range.Start
$fullText = $target.Range.Text
# Drop the trailing paragraph-mark character so the replacement range covers
# only the run text, leaving the paragraph mark itself untouched.
$textLen = $fullText.Length - 1

$old1 = "Before the workshop begins, you will receive a token that will grant you "
$old2 = "the access"
$old3 = " to the Git repository with code and related materials."

$expected = $old1 + $old2 + $old3
if ($fullText.Substring(0, $textLen) -ne $expected) {
    throw "Unexpected paragraph content: $fullText"
}

$editRange = $d.Range($pStart, $pStart + $textLen)

# Replace the three original runs (which were split around a gramStart/gramEnd
# proofing pair) with the two new runs called for by the edit, preserving the
# original run-level character formatting (sz/szCs 28) and the rsidRPr on the
# first run, and dropping the now-irrelevant grammar proofErr markers.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:r w:rsidRPr="00541B28"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Shortly before </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>the workshop, you will receive a token that will grant you access to the Git repository with code and related materials.</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$editRange.InsertXML($xml)

Write-Output ("Updated paragraph now reads: " + $target.Range.Text)
